$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 11:52"

# Row 4: Madrid
$ws.Range("B4").Value = 49526
$ws.Range("C4").Value = 27433
$ws.Range("D4").Value = 15369
$ws.Range("E4").Value = 6724

# Row 5: Cataluña
$ws.Range("C5").Value = 16651
$ws.Range("D5").Value = 16098

# Row 10: Navarra
$ws.Range("B10").Value = 4246
$ws.Range("C10").Value = 808
$ws.Range("D10").Value = 3186
$ws.Range("E10").Value = 252

# Row 12: La Rioja
$ws.Range("B12").Value = 3457
$ws.Range("C12").Value = 1585
$ws.Range("D12").Value = 1626
$ws.Range("E12").Value = 246

# Row 28: Cantabria
$ws.Range("B28").Value = 1823
$ws.Range("C28").Value = 363
$ws.Range("D28").Value = 1328
$ws.Range("E28").Value = 132

# Row 55: Ceuta
$ws.Range("B55").Value = 98
$ws.Range("C55").Value = 30
$ws.Range("D55").Value = 64
